$wb = $excel.ActiveWorkbook

# --- 1) Update existing values on "BBDD - Error Actual" sheet ---
$ws1 = $wb.Worksheets.Item("BBDD - Error Actual")

$ws1.Range("DF5").Value = -821169.4864866072
$ws1.Range("DG5").Value = -246350.8459459821
$ws1.Range("DH5").Value = -246350.8459459821
$ws1.Range("DI5").Value = -246350.8459459821
$ws1.Range("DJ5").Value = -82116.94864866079
$ws1.Range("DK5").Value = 739052.5378379466
$ws1.Range("DL5").Value = -223116.2690201853
$ws1.Range("DM5").Value = -217521.3600701397
$ws1.Range("DN5").Value = -213842.8122333892
$ws1.Range("DO5").Value = -46859.08442443195
$ws1.Range("EA5").Value = -223116.2690201853
$ws1.Range("EB5").Value = -440637.629090325
$ws1.Range("EC5").Value = -654480.4413237141
$ws1.Range("ED5").Value = -701339.525748146
$ws1.Range("EE5").Value = -663884.8863905093
$ws1.Range("EF5").Value = -624596.269608216
$ws1.Range("EG5").Value = -583740.8487692167
$ws1.Range("EH5").Value = -541514.7189122281
$ws1.Range("EI5").Value = -498070.3052207374
$ws1.Range("EJ5").Value = -453526.6418923203
$ws1.Range("EK5").Value = -407986.5396566967
$ws1.Range("EL5").Value = -361533.0312019755
$ws1.Range("EM5").Value = -314236.369219674
$ws1.Range("EN5").Value = -266154.697202705
$ws1.Range("EO5").Value = -217714.6452753017
$ws1.Range("FT5").Value = 69440.67433944368
$ws1.Range("DF6").Value = -172793.5867078734
$ws1.Range("DG6").Value = -51838.07601236203
$ws1.Range("DH6").Value = -51838.07601236203
$ws1.Range("DI6").Value = -51838.07601236203
$ws1.Range("DJ6").Value = -17279.35867078736
$ws1.Range("DK6").Value = 155514.2280370861
$ws1.Range("DL6").Value = -2694.949339040359
$ws1.Range("DM6").Value = -201.3831744067356
$ws1.Range("DN6").Value = 1930.179596015514
$ws1.Range("DO6").Value = 38354.50489018309
$ws1.Range("EA6").Value = -2694.949339040359
$ws1.Range("EB6").Value = -2896.332513447094
$ws1.Range("EC6").Value = -966.1529174315801
$ws1.Range("ED6").Value = 37388.35197275151
$ws1.Range("EP6").Value = 0
$ws1.Range("EQ6").Value = 0
$ws1.Range("ER6").Value = 0
$ws1.Range("ES6").Value = -10094.85503264291
$ws1.Range("FE6").Value = 49143.12667332167
$ws1.Range("FF6").Value = 51636.69283795529
$ws1.Range("FG6").Value = 53768.25560837754
$ws1.Range("FH6").Value = 45539.00852832754
$ws1.Range("FT6").Value = 368905.2651620624
$ws1.Range("DF7").Value = -172793.5867078734
$ws1.Range("DG7").Value = -51838.07601236203
$ws1.Range("DH7").Value = -51838.07601236203
$ws1.Range("DI7").Value = -51838.07601236203
$ws1.Range("DJ7").Value = -17279.35867078736
$ws1.Range("DK7").Value = 155514.2280370861
$ws1.Range("DL7").Value = -30222.69725944728
$ws1.Range("DM7").Value = -28950.38169411987
$ws1.Range("DN7").Value = -27880.67385196586
$ws1.Range("DO7").Value = 7602.90877774647
$ws1.Range("EA7").Value = -30222.69725944728
$ws1.Range("EB7").Value = -59173.07895356715
$ws1.Range("EC7").Value = -87053.752805533
$ws1.Range("ED7").Value = -79450.84402778653
$ws1.Range("EE7").Value = -53753.96635305126
$ws1.Range("EF7").Value = -27327.33242092184
$ws1.Range("EG7").Value = -243.0431073451546
$ws1.Range("EH7").Value = 27441.19741922177
$ws1.Range("EV7").Value = 0
$ws1.Range("EW7").Value = -7409.123303189877
$ws1.Range("FK7").Value = 27084.28931357669
$ws1.Range("FL7").Value = 20275.11722337704
$ws1.Range("FT7").Value = 151186.8788804009
$ws1.Range("DF8").Value = -172793.5867078734
$ws1.Range("DG8").Value = -51838.07601236203
$ws1.Range("DH8").Value = -51838.07601236203
$ws1.Range("DI8").Value = -51838.07601236203
$ws1.Range("DJ8").Value = -17279.35867078736
$ws1.Range("DK8").Value = 155514.2280370861
$ws1.Range("DL8").Value = -51028.21732314814
$ws1.Range("DM8").Value = -50982.54026582958
$ws1.Range("DN8").Value = -50943.94780441035
$ws1.Range("DO8").Value = -16351.74382640163
$ws1.Range("EA8").Value = -51028.21732314814
$ws1.Range("EB8").Value = -102010.7575889777
$ws1.Range("EC8").Value = -152954.7053933881
$ws1.Range("ED8").Value = -169306.4492197897
$ws1.Range("EE8").Value = -168349.259055316
$ws1.Range("EF8").Value = -167365.5171171398
$ws1.Range("EG8").Value = -166357.8050232289
$ws1.Range("EH8").Value = -165328.1945495052
$ws1.Range("EI8").Value = -164278.4276564388
$ws1.Range("EJ8").Value = -163209.9401613229
$ws1.Range("EK8").Value = -162125.5433885352
$ws1.Range("EL8").Value = -161041.1466157476
$ws1.Range("EM8").Value = -159956.7498429599
$ws1.Range("EN8").Value = -158872.3530701722
$ws1.Range("EO8").Value = -157787.9562973845
$ws1.Range("FT8").Value = -49650.30372026883

# --- 2) Add new worksheet "RESUMEN E1" at the end ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "RESUMEN E1"

# --- 3) Populate RESUMEN E1 sheet content ---
$ws4.Range("A1").Value = "Año"
$arrRow1 = New-Object 'object[,]' 1,15
$arrRow1[0,0] = 1
$arrRow1[0,1] = 2
$arrRow1[0,2] = 3
$arrRow1[0,3] = 4
$arrRow1[0,4] = 5
$arrRow1[0,5] = 6
$arrRow1[0,6] = 7
$arrRow1[0,7] = 8
$arrRow1[0,8] = 9
$arrRow1[0,9] = 10
$arrRow1[0,10] = 11
$arrRow1[0,11] = 12
$arrRow1[0,12] = 13
$arrRow1[0,13] = 14
$arrRow1[0,14] = 15
$ws4.Range("B1:P1").Value = $arrRow1

$ws4.Range("A2").Value = "Scenario 1: medidores con VAN > 0"

$ws4.Range("A3").Value = "Consumo (m3/año)"
$arrRow3 = New-Object 'object[,]' 1,15
$arrRow3[0,0] = 1608
$arrRow3[0,1] = 1608
$arrRow3[0,2] = 1608
$arrRow3[0,3] = 1608
$arrRow3[0,4] = 1608
$arrRow3[0,5] = 1608
$arrRow3[0,6] = 1608
$arrRow3[0,7] = 1608
$arrRow3[0,8] = 1608
$arrRow3[0,9] = 1608
$arrRow3[0,10] = 1608
$arrRow3[0,11] = 1608
$arrRow3[0,12] = 1608
$arrRow3[0,13] = 1608
$arrRow3[0,14] = 1608
$ws4.Range("B3:P3").Value = $arrRow3

$ws4.Range("A4").Value = "Submedición Con Proyecto (m3/año)"
$arrRow4 = New-Object 'object[,]' 1,15
$arrRow4[0,0] = 32.81632653061224
$arrRow4[0,1] = 32.81632653061224
$arrRow4[0,2] = 32.81632653061224
$arrRow4[0,3] = 32.81632653061224
$arrRow4[0,4] = 32.81632653061224
$arrRow4[0,5] = 32.81632653061224
$arrRow4[0,6] = 32.81632653061224
$arrRow4[0,7] = 32.81632653061224
$arrRow4[0,8] = 32.81632653061224
$arrRow4[0,9] = 32.81632653061224
$arrRow4[0,10] = 32.81632653061224
$arrRow4[0,11] = 32.81632653061224
$arrRow4[0,12] = 32.81632653061224
$arrRow4[0,13] = 32.81632653061224
$arrRow4[0,14] = 32.81632653061224
$ws4.Range("B4:P4").Value = $arrRow4

$ws4.Range("A5").Value = "Submedición Sin Proyecto (m3)"
$arrRow5 = New-Object 'object[,]' 1,15
$arrRow5[0,0] = 88.30886015065067
$arrRow5[0,1] = 93.8353731505637
$arrRow5[0,2] = 97.89714607092765
$arrRow5[0,3] = 101.1680835117484
$arrRow5[0,4] = 103.9252635949557
$arrRow5[0,5] = 106.3227143668892
$arrRow5[0,6] = 108.4373762139851
$arrRow5[0,7] = 110.3355575326812
$arrRow5[0,8] = 112.0587771467003
$arrRow5[0,9] = 113.6415221869177
$arrRow5[0,10] = 114.7840268451239
$arrRow5[0,11] = 115.3462306069298
$arrRow5[0,12] = 115.844019638902
$arrRow5[0,13] = 116.3074813602413
$arrRow5[0,14] = 116.519065114395
$ws4.Range("B5:P5").Value = $arrRow5

$ws4.Range("A6").Value = "Diferencial (volumen recuperado, m3/año)"
$arrRow6 = New-Object 'object[,]' 1,15
$arrRow6[0,0] = 55.49253362003843
$arrRow6[0,1] = 61.01904661995145
$arrRow6[0,2] = 65.0808195403154
$arrRow6[0,3] = 68.35175698113621
$arrRow6[0,4] = 71.10893706434345
$arrRow6[0,5] = 73.50638783627694
$arrRow6[0,6] = 75.62104968337287
$arrRow6[0,7] = 77.51923100206899
$arrRow6[0,8] = 79.24245061608802
$arrRow6[0,9] = 80.82519565630548
$arrRow6[0,10] = 81.96770031451162
$arrRow6[0,11] = 82.52990407631758
$arrRow6[0,12] = 83.0276931082898
$arrRow6[0,13] = 83.49115482962902
$arrRow6[0,14] = 83.70273858378276
$ws4.Range("B6:P6").Value = $arrRow6

$ws4.Range("A7").Value = "Diferencial (% c/r consumo renovados)"
$arrRow7 = New-Object 'object[,]' 1,15
$arrRow7[0,0] = 0.03451028210201395
$arrRow7[0,1] = 0.03794716829598971
$arrRow7[0,2] = 0.04047314648029565
$arrRow7[0,3] = 0.04250731155543296
$arrRow7[0,4] = 0.04422197578628324
$arrRow7[0,5] = 0.0457129277588787
$arrRow7[0,6] = 0.04702801597224681
$arrRow7[0,7] = 0.04820847699133644
$arrRow7[0,8] = 0.04928013098015424
$arrRow7[0,9] = 0.05026442515939396
$arrRow7[0,10] = 0.05097493800653707
$arrRow7[0,11] = 0.05132456721164028
$arrRow7[0,12] = 0.05163413750515534
$arrRow7[0,13] = 0.05192235996867477
$arrRow7[0,14] = 0.05205394190533754
$ws4.Range("B7:P7").Value = $arrRow7

$ws4.Range("A8").Value = "Ingresos volumen recuperado ($/año)"
$arrRow8 = New-Object 'object[,]' 1,15
$arrRow8[0,0] = 97655.15049561902
$arrRow8[0,1] = 107380.6473060155
$arrRow8[0,2] = 114528.5106299914
$arrRow8[0,3] = 120284.6703727097
$arrRow8[0,4] = 125136.7255665259
$arrRow8[0,5] = 129355.7330737712
$arrRow8[0,6] = 133077.0917404969
$arrRow8[0,7] = 136417.4903536592
$arrRow8[0,8] = 139449.9932827248
$arrRow8[0,9] = 142235.2906014043
$arrRow8[0,10] = 144245.8577364896
$arrRow8[0,11] = 145235.2177347965
$arrRow8[0,12] = 146111.221399815
$arrRow8[0,13] = 146926.816241028
$arrRow8[0,14] = 147299.1590051156
$ws4.Range("B8:P8").Value = $arrRow8

$ws4.Range("A9").Value = "Con Proyecto - Error ponderado final"
$arrRow9 = New-Object 'object[,]' 1,15
$arrRow9[0,0] = -0.02
$arrRow9[0,1] = -0.02
$arrRow9[0,2] = -0.02
$arrRow9[0,3] = -0.02
$arrRow9[0,4] = -0.02
$arrRow9[0,5] = -0.02
$arrRow9[0,6] = -0.02
$arrRow9[0,7] = -0.02
$arrRow9[0,8] = -0.02
$arrRow9[0,9] = -0.02
$arrRow9[0,10] = -0.02
$arrRow9[0,11] = -0.02
$arrRow9[0,12] = -0.02
$arrRow9[0,13] = -0.02
$arrRow9[0,14] = -0.02
$ws4.Range("B9:P9").Value = $arrRow9

$ws4.Range("A10").Value = "Sin Proyecto - Error ponderado final"
$arrRow10 = New-Object 'object[,]' 1,15
$arrRow10[0,0] = -0.05205942280040198
$arrRow10[0,1] = -0.05513774988520111
$arrRow10[0,2] = -0.05738748452473071
$arrRow10[0,3] = -0.05919141861336603
$arrRow10[0,4] = -0.06070665922454929
$arrRow10[0,5] = -0.0620202447741322
$arrRow10[0,6] = -0.06317584184351066
$arrRow10[0,7] = -0.06421071661411089
$arrRow10[0,8] = -0.06514822553482019
$arrRow10[0,9] = -0.06600765648505295
$arrRow10[0,10] = -0.06662705542686274
$arrRow10[0,11] = -0.06693154779832435
$arrRow10[0,12] = -0.06720098704937827
$arrRow10[0,13] = -0.06745170604287506
$arrRow10[0,14] = -0.06756612175039409
$ws4.Range("B10:P10").Value = $arrRow10

# --- 4) Apply header style (row 1) matching existing bold/border/center style ---
# Reuse the style already present on row 1 of the first sheet (bold, centered, thin border)
# by copying its format only, so we do not introduce extra unused style entries.
$styleSource = $ws1.Range("A1:P1")
$styleSource.Copy()
$ws4.Range("A1:P1").PasteSpecial(-4122)

# --- 5) Restore original active sheet/tab selection ---
$ws1.Activate()

Write-Host "Edit applied successfully"